$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "UnitMass" column (C) values in the first table (rows 2-21)
$ws.Range("C2").Value = 39
$ws.Range("C3").Value = 23
$ws.Range("C4").Value = 57
$ws.Range("C5").Value = 113
$ws.Range("C6").Value = 59
$ws.Range("C7").Value = 175
$ws.Range("C8").Value = 69
$ws.Range("C9").Value = 231
$ws.Range("C10").Value = 112
$ws.Range("C11").Value = 115
$ws.Range("C12").Value = 111
$ws.Range("C13").Value = 96
$ws.Range("C14").Value = 159
$ws.Range("C15").Value = 99
$ws.Range("C16").Value = 83
$ws.Range("C17").Value = 82
$ws.Range("C18").Value = 103
$ws.Range("C19").Value = 119
$ws.Range("C20").Value = 81
$ws.Range("C21").Value = 98

# Update "UnitMass" column (C) values in the second table (rows 23-42)
$ws.Range("C23").Value = 27
$ws.Range("C24").Value = 28
$ws.Range("C25").Value = 46
$ws.Range("C26").Value = 71
$ws.Range("C27").Value = 91
$ws.Range("C28").Value = 102
$ws.Range("C29").Value = 88
$ws.Range("C30").Value = 74
$ws.Range("C31").Value = 44
$ws.Range("C32").Value = 45
$ws.Range("C33").Value = 29
$ws.Range("C34").Value = 56
$ws.Range("C35").Value = 77
$ws.Range("C36").Value = 1
$ws.Range("C37").Value = 30
$ws.Range("C38").Value = 149
$ws.Range("C39").Value = 61
$ws.Range("C40").Value = 51
$ws.Range("C41").Value = 73
$ws.Range("C42").Value = 24
